# Generate Report for Handback
# - Updates status text from "Ready for handoff" to "Handed back: in sync with en-US"
#   on the Overview sheet and on each language sheet's Status column.
# - Fills in "Latest Target File" (hyperlink to the source .md) and
#   "Latest Handback File" (generated .xlf file name) columns for both rows
#   on the zh-cn and de-de sheets.
# - Updates "Latest Handback DateTime" with the new handback timestamps.
# - Widens a few columns that now hold longer text.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"
$hyperlinkColor = 15570276   # BGR for RGB(100,149,237) / #6495ED, matches workbook's HyperLink style

$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d8ceea8c8a685b63bb59ad6754812cf729d318d0/e2e/4a066e60-9c9f-4238-a567-d31be1f8a179.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d8ceea8c8a685b63bb59ad6754812cf729d318d0/e2e/523eb56d-c432-4139-b454-edeb6f0553e4.md"
$mdName1 = "4a066e60-9c9f-4238-a567-d31be1f8a179.md"
$mdName2 = "523eb56d-c432-4139-b454-edeb6f0553e4.md"

# ---------------------------------------------------------------------------
# Overview sheet: refresh the per-language status cells
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText
$wsOverview.Columns.Item(5).ColumnWidth = 29.14
$wsOverview.Columns.Item(6).ColumnWidth = 29.14

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdUrl1, "", "", $mdName1)
$wsZh.Range("I2").Font.Color = $hyperlinkColor
$wsZh.Range("I2").Font.Underline = $true
$wsZh.Range("J2").Value = "4a066e60-9c9f-4238-a567-d31be1f8a179.df1b08eeb98abb1bdea2c68622bef2fdf8082352.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-10-10 09:37:27"

$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $mdUrl2, "", "", $mdName2)
$wsZh.Range("I3").Font.Color = $hyperlinkColor
$wsZh.Range("I3").Font.Underline = $true
$wsZh.Range("J3").Value = "523eb56d-c432-4139-b454-edeb6f0553e4.4d2b0affa740e69a3be9f3db847d366d8bbe0601.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-10-10 09:37:27"

$wsZh.Columns.Item(3).ColumnWidth = 29.14
$wsZh.Columns.Item(9).ColumnWidth = 39.14
$wsZh.Columns.Item(10).ColumnWidth = 39.14

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdUrl1, "", "", $mdName1)
$wsDe.Range("I2").Font.Color = $hyperlinkColor
$wsDe.Range("I2").Font.Underline = $true
$wsDe.Range("J2").Value = "4a066e60-9c9f-4238-a567-d31be1f8a179.df1b08eeb98abb1bdea2c68622bef2fdf8082352.de-de.xlf"
$wsDe.Range("K2").Value = "2016-10-10 09:37:43"

$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $mdUrl2, "", "", $mdName2)
$wsDe.Range("I3").Font.Color = $hyperlinkColor
$wsDe.Range("I3").Font.Underline = $true
$wsDe.Range("J3").Value = "523eb56d-c432-4139-b454-edeb6f0553e4.4d2b0affa740e69a3be9f3db847d366d8bbe0601.de-de.xlf"
$wsDe.Range("K3").Value = "2016-10-10 09:37:43"

$wsDe.Columns.Item(3).ColumnWidth = 29.14
$wsDe.Columns.Item(9).ColumnWidth = 39.14
$wsDe.Columns.Item(10).ColumnWidth = 39.14

Write-Output "Handback report generated"
